$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update "last updated" timestamp (row 1)
$ws.Range("A1").Value = "Datos actualizados a 26 de Julio de 2020 a las 23:30"

# --- Updated country statistics (Casos totales, Nuevos casos, Casos activos,
#     Recuperados, Casos criticos, Muertes hoy, Muertes) ---
# Columns: A=Pais B=Casos totales C=Nuevos casos D=Casos activos
#          E=Recuperados F=Casos criticos G=Muertes hoy H=Muertes

# Row 4: Estados Unidos
$ws.Range("B4").Value = 4362291
$ws.Range("C4").Value = 46582
$ws.Range("D4").Value = 2084599
$ws.Range("E4").Value = 2127952
$ws.Range("F4").Value = 0
$ws.Range("G4").Value = 342
$ws.Range("H4").Value = 149740

# Row 5: Brasil
$ws.Range("B5").Value = 2419091
$ws.Range("C5").Value = 22657
$ws.Range("D5").Value = 1617480
$ws.Range("E5").Value = 714607
$ws.Range("F5").Value = 0
$ws.Range("G5").Value = 508
$ws.Range("H5").Value = 87004

# Row 21: Alemania
$ws.Range("B21").Value = 206741
$ws.Range("C21").Value = 409
$ws.Range("D21").Value = 190600
$ws.Range("E21").Value = 6938
$ws.Range("F21").Value = 0
$ws.Range("G21").Value = 1
$ws.Range("H21").Value = 9203

# Row 28: Egipto
$ws.Range("B28").Value = 92062
$ws.Range("C28").Value = 479
$ws.Range("D28").Value = 33831
$ws.Range("E28").Value = 53625
$ws.Range("F28").Value = 0
$ws.Range("G28").Value = 48
$ws.Range("H28").Value = 4606

# Row 41: Israel
$ws.Range("B41").Value = 61956
$ws.Range("C41").Value = 1278
$ws.Range("D41").Value = 27025
$ws.Range("E41").Value = 34461
$ws.Range("F41").Value = 0
$ws.Range("G41").Value = 13
$ws.Range("H41").Value = 470

# Row 70: Costa de Marfil
$ws.Range("B70").Value = 15596
$ws.Range("C70").Value = 102
$ws.Range("D70").Value = 10178
$ws.Range("E70").Value = 5322
$ws.Range("F70").Value = 0
$ws.Range("G70").Value = 2
$ws.Range("H70").Value = 96

# Row 80: Estado de Palestina
$ws.Range("B80").Value = 10469
$ws.Range("C80").Value = 163
$ws.Range("D80").Value = 3752
$ws.Range("E80").Value = 6641
$ws.Range("F80").Value = 0
$ws.Range("G80").Value = 1
$ws.Range("H80").Value = 76

# Rows 116-117: Libia moves above Hong Kong (reordered by Casos totales)
$ws.Range("A116").Value = "Libia"
$ws.Range("B116").Value = 2669
$ws.Range("C116").Value = 122
$ws.Range("D116").Value = 553
$ws.Range("E116").Value = 2056
$ws.Range("F116").Value = 0
$ws.Range("G116").Value = 2
$ws.Range("H116").Value = 60

$ws.Range("A117").Value = "Hong Kong"
$ws.Range("B117").Value = 2634
$ws.Range("C117").Value = 128
$ws.Range("D117").Value = 1495
$ws.Range("E117").Value = 1121
$ws.Range("F117").Value = 0
$ws.Range("G117").Value = 0
$ws.Range("H117").Value = 18

# Rows 130-134: Ruanda & Sierra Leona move above Namibia & Benin (reordered
# by Casos totales); Yemen keeps its row but gets updated figures.
$ws.Range("A130").Value = "Ruanda"
$ws.Range("B130").Value = 1821
$ws.Range("C130").Value = 69
$ws.Range("D130").Value = 918
$ws.Range("E130").Value = 898
$ws.Range("F130").Value = 0
$ws.Range("G130").Value = 0
$ws.Range("H130").Value = 5

$ws.Range("A131").Value = "Sierra Leona"
$ws.Range("B131").Value = 1783
$ws.Range("C131").Value = 15
$ws.Range("D131").Value = 1313
$ws.Range("E131").Value = 404
$ws.Range("F131").Value = 0
$ws.Range("G131").Value = 0
$ws.Range("H131").Value = 66

$ws.Range("A132").Value = "Namibia"
$ws.Range("B132").Value = 1775
$ws.Range("C132").Value = 88
$ws.Range("D132").Value = 75
$ws.Range("E132").Value = 1692
$ws.Range("F132").Value = 0
$ws.Range("G132").Value = 1
$ws.Range("H132").Value = 8

$ws.Range("A133").Value = "Benin"
$ws.Range("B133").Value = 1770
$ws.Range("C133").Value = 76
$ws.Range("D133").Value = 1036
$ws.Range("E133").Value = 699
$ws.Range("F133").Value = 0
$ws.Range("G133").Value = 1
$ws.Range("H133").Value = 35

# Row 134: Yemen
$ws.Range("B134").Value = 1681
$ws.Range("C134").Value = 7
$ws.Range("D134").Value = 797
$ws.Range("E134").Value = 405
$ws.Range("F134").Value = 0
$ws.Range("G134").Value = 5
$ws.Range("H134").Value = 479
